$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression: only B2 changes
$ws.Range("B2").Value = 7893751996506930

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 2942555796670.181
$ws.Range("C3").Value = 2891072235705.228
$ws.Range("D3").Value = 318569518135916.7

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 2868291126679.81
$ws.Range("C4").Value = 2826983785464.827
$ws.Range("D4").Value = 84324533490112.72

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 414592087959714.8
$ws.Range("C5").Value = 900214740170152.9
$ws.Range("D5").Value = 2657507705113638
